# Generate Report for handback
# Updates the zh-cn and de-de status sheets to reflect that the handback
# has occurred: status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Target File / Latest
# Handback File columns are populated with hyperlinks to the handed-back
# files, and the Latest Handback DateTime is stamped.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 2 (715bcfdd-b5b1-4701-b415-8cdb3465f613)
$ws.Range("B2").Value = $statusHandedBack
$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/37baf36f718dce052dde0d7053ceb6d1f405d519/e2e/715bcfdd-b5b1-4701-b415-8cdb3465f613.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "715bcfdd-b5b1-4701-b415-8cdb3465f613.md"
) | Out-Null
$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df490e2fc9e24b684af53eb18d62131c1f82204e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/715bcfdd-b5b1-4701-b415-8cdb3465f613.cf4212ffff762640215569d199e0f6e656b94972.zh-cn.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "715bcfdd-b5b1-4701-b415-8cdb3465f613.cf4212ffff762640215569d199e0f6e656b94972.zh-cn.xlf"
) | Out-Null
$ws.Range("G2").Value = "2016-01-25 14:12:40"

# Row 3 (bd62af68-9bda-4622-a462-1f5ffd07dcf5)
$ws.Range("B3").Value = $statusHandedBack
$ws.Hyperlinks.Add(
    $ws.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/37baf36f718dce052dde0d7053ceb6d1f405d519/e2e/bd62af68-9bda-4622-a462-1f5ffd07dcf5.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "bd62af68-9bda-4622-a462-1f5ffd07dcf5.md"
) | Out-Null
$ws.Hyperlinks.Add(
    $ws.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df490e2fc9e24b684af53eb18d62131c1f82204e/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/bd62af68-9bda-4622-a462-1f5ffd07dcf5.b3cf49591cb83c962e75ca44ca1a5a23a5b4c6db.zh-cn.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "bd62af68-9bda-4622-a462-1f5ffd07dcf5.b3cf49591cb83c962e75ca44ca1a5a23a5b4c6db.zh-cn.xlf"
) | Out-Null
$ws.Range("G3").Value = "2016-01-25 14:12:40"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 2 (715bcfdd-b5b1-4701-b415-8cdb3465f613)
$ws.Range("B2").Value = $statusHandedBack
$ws.Hyperlinks.Add(
    $ws.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/37baf36f718dce052dde0d7053ceb6d1f405d519/e2e/715bcfdd-b5b1-4701-b415-8cdb3465f613.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "715bcfdd-b5b1-4701-b415-8cdb3465f613.md"
) | Out-Null
$ws.Hyperlinks.Add(
    $ws.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a2d976bf14c6eb49d68b05b40fe8424e5e63a67/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/715bcfdd-b5b1-4701-b415-8cdb3465f613.cf4212ffff762640215569d199e0f6e656b94972.de-de.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "715bcfdd-b5b1-4701-b415-8cdb3465f613.cf4212ffff762640215569d199e0f6e656b94972.de-de.xlf"
) | Out-Null
$ws.Range("G2").Value = "2016-01-25 14:12:59"

# Row 3 (bd62af68-9bda-4622-a462-1f5ffd07dcf5)
$ws.Range("B3").Value = $statusHandedBack
$ws.Hyperlinks.Add(
    $ws.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/37baf36f718dce052dde0d7053ceb6d1f405d519/e2e/bd62af68-9bda-4622-a462-1f5ffd07dcf5.md",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "bd62af68-9bda-4622-a462-1f5ffd07dcf5.md"
) | Out-Null
$ws.Hyperlinks.Add(
    $ws.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a2d976bf14c6eb49d68b05b40fe8424e5e63a67/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/bd62af68-9bda-4622-a462-1f5ffd07dcf5.b3cf49591cb83c962e75ca44ca1a5a23a5b4c6db.de-de.xlf",
    [System.Type]::Missing,
    [System.Type]::Missing,
    "bd62af68-9bda-4622-a462-1f5ffd07dcf5.b3cf49591cb83c962e75ca44ca1a5a23a5b4c6db.de-de.xlf"
) | Out-Null
$ws.Range("G3").Value = "2016-01-25 14:12:59"

Write-Host "Handback report generated."
